$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: GDP Nowcast ---
# Date moves forward one quarter and the highlight ("as of" fill) moves onto C7.
$ws.Range("C7").Value = 45931
$ws.Range("C3").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rolling 5-period history shifts right (oldest dropped, new reading prepended).
$ws.Range("F7").Value = 2.9721
$ws.Range("G7").Value = 3.4728
$ws.Range("H7").Value = 2.902
$ws.Range("I7").Value = -2.7318
$ws.Range("J7").Value = 2.2711

# --- Rows 18-21: CPI / Core CPI blocks lose their highlight (moved elsewhere) ---
$ws.Range("N3").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N3").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N3").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N3").Copy()
$ws.Range("N21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 29: 5yr, 5yr Forward breakeven ---
$ws.Range("N29").Value = 46015
$ws.Range("R29").Value = 2.24
$ws.Range("T29").Value = 2.21

# --- Row 30: 10yr TIPS breakeven ---
$ws.Range("N30").Value = 46015
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 2.23

# --- Row 51: 30y Mortgage rate ---
$ws.Range("N51").Value = 46013
$ws.Range("Q51").Value = 6.18
$ws.Range("R51").Value = 6.21
$ws.Range("S51").Value = 6.22
$ws.Range("T51").Value = 6.19
$ws.Range("U51").Value = 6.23

# --- Row 52: BAA corporate bond yield ---
$ws.Range("N52").Value = 46014
$ws.Range("Q52").Value = 5.92
$ws.Range("R52").Value = 5.93
$ws.Range("S52").Value = 5.92
$ws.Range("T52").Value = 5.9
$ws.Range("U52").Value = 5.94
